$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 783
$ws1.Range("F5").Value = 843
$ws1.Range("F6").Value = 2083
$ws1.Range("F7").Value = 187

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 783
$ws4.Range("F7").Value = 843
$ws4.Range("F8").Value = 2083
$ws4.Range("F10").Value = 187
